# Generate Report for Handoff
# Re-generate the localization-status report: the 3203bb0e-... entry has
# been handed off (status -> "Ready for handoff" with new timestamps) and
# moved to the bottom of each table, while the other two rows shift up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-21 12:56:36"

$ws.Range("A3").Value = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-21 12:56:36"

$ws.Range("A4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-21 12:58:49"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md" }
    elseif ($addr -eq '$A$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md" }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-21 12:56:32"
$ws.Range("F2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$ws.Range("G2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-21 12:56:54"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 12:56:32"
$ws.Range("F3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$ws.Range("G3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-21 12:56:54"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-21 12:58:45"
$ws.Range("F4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md"
$ws.Range("G4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-21 12:58:12"
$ws.Range("J4").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md" }
    elseif ($addr -eq '$D$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf" }
    elseif ($addr -eq '$F$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.md" }
    elseif ($addr -eq '$G$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md" }
    elseif ($addr -eq '$D$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf" }
    elseif ($addr -eq '$F$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.md" }
    elseif ($addr -eq '$G$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.zh-cn.xlf" }
    elseif ($addr -eq '$A$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md" }
    elseif ($addr -eq '$D$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.zh-cn.xlf" }
    elseif ($addr -eq '$F$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md" }
    elseif ($addr -eq '$G$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.zh-cn.xlf" }
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$ws.Range("E2").Value = "2016-03-21 12:56:36"
$ws.Range("F2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$ws.Range("G2").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$ws.Range("H2").Value = "2016-03-21 12:57:00"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 12:56:36"
$ws.Range("F3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.md"
$ws.Range("G3").Value = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf"
$ws.Range("H3").Value = "2016-03-21 12:57:00"
$ws.Range("J3").Value = "Include"

$ws.Range("A4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.de-de.xlf"
$ws.Range("E4").Value = "2016-03-21 12:58:49"
$ws.Range("F4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md"
$ws.Range("G4").Value = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.de-de.xlf"
$ws.Range("H4").Value = "2016-03-21 12:58:17"
$ws.Range("J4").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff8d28186a-9f8f-4d6e-af74-3eea744c3452.md" }
    elseif ($addr -eq '$D$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf" }
    elseif ($addr -eq '$F$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.md" }
    elseif ($addr -eq '$G$2') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "ffffff22af7fe0-cb28-4f36-8b19-b72962c87b02.md" }
    elseif ($addr -eq '$D$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf" }
    elseif ($addr -eq '$F$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.md" }
    elseif ($addr -eq '$G$3') { $h.TextToDisplay = "e11abc11-d81e-48f3-a5af-a03605f53f72.e5c2defb3d53fad5e7a5b8091116ba601f6445fc.de-de.xlf" }
    elseif ($addr -eq '$A$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md" }
    elseif ($addr -eq '$D$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.de-de.xlf" }
    elseif ($addr -eq '$F$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.md" }
    elseif ($addr -eq '$G$4') { $h.TextToDisplay = "3203bb0e-9de9-4904-8aa9-9eb6ec101252.bedd3d7b00d6cf81d85422e854f594792553aa8b.de-de.xlf" }
}

$wb.Save()
